$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$forceTextCells = @("D5", "D6", "D10", "D11", "D18", "D20", "D21", "D25", "D27", "D28", "D29", "D31", "D32", "D36", "D38", "D41", "D42", "D46", "D47", "D50")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.985.88"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "214.03"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").Value = "18.52"
$ws.Range("E10").Value = "  -5.70%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "1.860.43"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "1.633.81"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D16").Value = "26.000.28"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "0.0₃0745"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").Value = "61.77"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "190.11"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").Value = "4.24"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "143.21"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "1.76"
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "6.77"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("D29").Value = "15.19"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").Value = "3.15"
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "0.871"
$ws.Range("E36").Value = "  -3.74%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.135.42"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "2.43"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "98.61"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "0.781"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("D44").Value = "1.770.95"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("D46").Value = "55.11"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").Value = "0.0530"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D50").Value = "7.54"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("E51").Value = "  +0.21%  "

foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
